$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 493, shifting rows 493:596 down to 494:597
$ws.Rows.Item(493).Insert()

# Populate the new row 493 with the inserted record's data
$ws.Cells.Item(493, 1).Value = 3
$ws.Cells.Item(493, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(493, 3).Value = "Coquimbo"
$ws.Cells.Item(493, 4).Value = 44995
$ws.Cells.Item(493, 5).Value = 5
$ws.Cells.Item(493, 6).Value = "Fruta"
$ws.Cells.Item(493, 7).Value = 100108
$ws.Cells.Item(493, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(493, 9).Value = 100108002
$ws.Cells.Item(493, 10).Value = "Mango"
$ws.Cells.Item(493, 11).Value = "Sin especificar"
$ws.Cells.Item(493, 12).Value = "Primera"
$ws.Cells.Item(493, 13).Value = 456
$ws.Cells.Item(493, 14).Value = 7000
$ws.Cells.Item(493, 15).Value = 7000
$ws.Cells.Item(493, 16).Value = 7000
$ws.Cells.Item(493, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(493, 18).Value = "Perú"
$ws.Cells.Item(493, 19).Value = 1750
$ws.Cells.Item(493, 20).Value = 4
